$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 137 (shifts existing rows 137:156 down to 138:157)
$ws.Rows.Item(137).Insert()

# Populate the newly inserted row 137 with the new weekly record
$ws.Cells.Item(137, 1).Value  = 8
$ws.Cells.Item(137, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(137, 3).Value  = "Coquimbo"
$ws.Cells.Item(137, 4).Value  = 45127
$ws.Cells.Item(137, 5).Value  = 4
$ws.Cells.Item(137, 6).Value  = 100114007
$ws.Cells.Item(137, 7).Value  = "Jengibre"
$ws.Cells.Item(137, 8).Value  = "Sin especificar"
$ws.Cells.Item(137, 9).Value  = "Primera"
$ws.Cells.Item(137, 10).Value = 400
$ws.Cells.Item(137, 11).Value = 17500
$ws.Cells.Item(137, 12).Value = 18000
$ws.Cells.Item(137, 13).Value = 17750
$ws.Cells.Item(137, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(137, 15).Value = "Perú"
$ws.Cells.Item(137, 16).Value = 1365
$ws.Cells.Item(137, 17).Value = 13
$ws.Cells.Item(137, 18).Value = "Hortaliza"
